# Fixed start_with embedded function to be more flexible on arguments
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Row 1 (header) - re-point a few header cells
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "STATION, BEARINGS"
$ws.Range("C1").Value = "STATION, BEARINGS,ASSIGNMENT_NUMBER"
$ws.Range("D1").Value = "AAR_CAR_TYPE,CAR_SERIES,SCS"
$ws.Range("E1").Value = "AAR_CAR_TYPE,CAR_SERIES,SCS"

# ---------------------------------------------------------------------------
# 2) Row 2
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = "STATION > 6300"
$ws.Range("C2").Value = "STATION < 9714"
$ws.Range("D2").Value = 'starts_with(AAR_CAR_TYPE, "V") '
$ws.Range("E2").Value = 'starts_with(AAR_CAR_TYPE, "V") '

# ---------------------------------------------------------------------------
# 3) Row 3
# ---------------------------------------------------------------------------
$ws.Range("B3").Value = "STATION > 7330"
$ws.Range("C3").Value = "STATION < 9999"
$ws.Range("D3").Value = 'starts_with(SCS, "112J")'
$ws.Range("E3").Value = 'starts_with(SCS, "112J")'

# ---------------------------------------------------------------------------
# 4) Row 4
# ---------------------------------------------------------------------------
$ws.Range("D4").Value = 'exclude(CAR_SERIES ,"MILW" )>= 120000 '
$ws.Range("E4").Value = 'exclude(CAR_SERIES ,"MILW") <= 120209'

# ---------------------------------------------------------------------------
# 5) Row 5 - D5 gets a new formula value
# ---------------------------------------------------------------------------
$ws.Range("D5").Value = " AAR_CAR_TYPE = ['M310','M340']"

# ---------------------------------------------------------------------------
# 6) Row 6 - rule list, bearings, formulas and combined error message
# ---------------------------------------------------------------------------
$ws.Range("A6").Value = "Rule 41,46,42,47"
$ws.Range("B6").Value = "BEARINGS = ['A','C']"
$ws.Range("D6").Value = ' starts_with(AAR_CAR_TYPE, "T")'
$ws.Range("E6").Value = " "
$ws.Range("H6").Value = "******************************************************`nTRAIN HANDLING TANK CAR(S) WITH PLAIN BEARINGS OR`nROLLER BEARINGS WITH CONVERTED FRICTION BEARING`nTRUCK SIDEFRAMES.   CARS ARE PROHIBITED IN`nINTERCHANGE OR MOVEMENT ON CPR IF: CONTAINING OR LAST CONTAINED DANGEROUS GOODS`n******************************************************"
$ws.Rows.Item(6).RowHeight = 182.25

# ---------------------------------------------------------------------------
# 7) Row 7 - rule 57 / assignment number / clear B7 / new UTLX message
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = "Rule 57"
$ws.Range("B7").ClearContents()
$ws.Range("C7").Value = "ASSIGNMENT_NUMBER=5000"
$ws.Range("D7").Value = " "
$ws.Range("H7").Value = "******************************************************`nTRAIN HANDLING TANK CARS IN UTLX SERIES.  CARS CANNOT`nMOVE BEYOND NEXT MECHANICAL INSPECTION FACILITY`n******************************************************"
$ws.Rows.Item(7).RowHeight = 122.25

# ---------------------------------------------------------------------------
# 8) New rows 8-11: copy formatting pattern from row 7, then fill in values
# ---------------------------------------------------------------------------
$ws.Range("A7:S7").Copy()
$ws.Range("A8:S8").PasteSpecial(-4122)
$ws.Range("A7:S7").Copy()
$ws.Range("A9:S9").PasteSpecial(-4122)
$ws.Range("A7:S7").Copy()
$ws.Range("A10:S10").PasteSpecial(-4122)
$ws.Range("A7:S7").Copy()
$ws.Range("A11:S11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

for ($r = 8; $r -le 11; $r++) {
  foreach ($col in @("I","J","K","L","M","N","O","P","Q","R","S")) {
    $ws.Range($col + $r).Value = "None"
  }
}

# Row 8 - Rule 59 / business-car data
$ws.Range("A8").Value = "Rule 59"
$ws.Range("B8").Value = ""
$ws.Range("C8").Value = ""
$ws.Range("D8").Value = " AAR_CAR_TYPE = 'M530'"
$ws.Range("E8").Value = "CAR_SERIES = ['CP  000070','CP  000085', 'CP  000095', 'CP  000096', 'CP  029114', 'CP  401750','CP  401753',  'CP  000099','CP  000102', 'CP  000105', 'CP  000106', 'CP  000110' ]"
$ws.Range("H8").Value = "******************************************************`nTRAIN HANDLING BUSINESS CAR(S).`nMUST BE MARSHALLED AS PER GOI SEC 7 ITEM 21.2`nB END TRAILING WHEN PRACTICABLE`nDO NOT EXCEED TIME TABLE FREIGHT TRAIN SPEED  OR AS INDICATED IN GOI SEC 7 ITEM 21.1`n******************************************************"
$ws.Rows.Item(8).RowHeight = 137.25

# Row 9 - Rule 178,179 / SOO fuel cars
$ws.Range("A9").Value = "Rule 178,179"
$ws.Range("C9").Value = "ASSIGNMENT_NUMBER=2240"
$ws.Range("D9").Value = " "
$ws.Range("E9").Value = " "
$ws.Range("F9").Value = "EMPTY_LOAD = 2"
$ws.Range("H9").Value = "******************************************************`nSOO 4000 4001 4002 AND 4003`nARE USED TO FUEL LOCOMOTIVES ENROUTE.`nCARS ARE EXEMPT FROM THE`nTRANSPORTATION OF DANGEROUS GOODS REGULATIONS  WHEN MARSHALLED NEXT TO LOCOMOTIVES`n******************************************************"
$ws.Rows.Item(9).RowHeight = 152.25

# Row 10 - Rule 56 / Detroit tunnel
$ws.Range("A10").Value = "Rule 56"
$ws.Range("B10").Value = "STATION = [100, 4544]"
$ws.Range("E10").Value = " "
$ws.Range("F10").Value = "EMPTY_LOAD = 2"
$ws.Range("H10").Value = "******************************************************`nTRAIN HANDLING MULTI LEVEL CAR(S) PROHIBITED`nIN BOTH TUBES OF THE DETROIT TUNNEL`n******************************************************"
$ws.Rows.Item(10).RowHeight = 92.25

# Row 11 - Rule 250,252 / articulate multi-platform cars
$ws.Range("A11").Value = "Rule 250,252"
$ws.Range("B11").Value = "STATION = [100, 4544]"
$ws.Range("C11").Value = "ASSIGNMENT_NUMBER=2240"
$ws.Range("D11").Value = 'starts_with(AAR_CAR_TYPE ,"S*6" )'
$ws.Range("E11").Value = ' starts_with(AAR_CAR_TYPE ,"S*8" )'
$ws.Range("F11").Value = "EMPTY_LOAD = 1"
$ws.Range("H11").Value = "******************************************************`nTRAIN HANDLING ARTICULATE MULTI-PLATFORM CAR(S)`nEQUIPPED WITH 125 TON TRUCKS`nLOADED WITH ONE OR MORE CONTAINERS`nIF CONTAINERS ARE DOUBLE STACKED SPEED RESTRICTIONS APPLY PER GOI SECTION 10`n******************************************************"
$ws.Rows.Item(11).RowHeight = 152.25

# ---------------------------------------------------------------------------
# 9) Column widths (best effort; COM rounds to whole-pixel granularity)
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 21.14
$ws.Columns.Item(3).ColumnWidth = 40.86
$ws.Columns.Item(9).ColumnWidth = 20.43

# ---------------------------------------------------------------------------
# 10) Sheet view: top-left cell + selection
# ---------------------------------------------------------------------------
$ws.Range("C9").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 2
